$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.417.69"
$ws.Range("E2").Value = "  +3.18%  "
$ws.Range("D3").Value = "1.582.95"
$ws.Range("E3").Value = "  +0.65%  "
$ws.Range("E4").Value = "  +1.13%  "
$ws.Range("D5").Value = "'212.35"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +0.59%  "
$ws.Range("D6").Value = "'0.492"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  -0.17%  "
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").Value = "'24.05"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.59%  "
$ws.Range("D9").Value = "'0.250"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  -0.13%  "
$ws.Range("D10").Value = "'0.0598"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "  +0.10%  "
$ws.Range("D11").Value = "'0.0883"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +1.18%  "
$ws.Range("D12").Value = "1.810.07"
$ws.Range("E12").Value = "  +0.73%  "
$ws.Range("D13").Value = "1.593.43"
$ws.Range("E13").Value = "  +1.28%  "
$ws.Range("D14").Value = "'0.527"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.06%  "
$ws.Range("D15").Value = "'3.72"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  -1.22%  "
$ws.Range("D16").Value = "28.444.55"
$ws.Range("E16").Value = "  +3.47%  "
$ws.Range("D17").Value = "'62.70"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  -0.87%  "
$ws.Range("D18").Value = "'229.43"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.39%  "
$ws.Range("D19").Value = "'7.42"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  -1.51%  "
$ws.Range("D20").Value = "0.0₃0701"
$ws.Range("E20").Value = "  -0.81%  "
$ws.Range("E21").Value = "  +1.00%  "
$ws.Range("D22").Value = "'3.99"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  -3.31%  "
$ws.Range("D23").Value = "'9.26"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  -2.09%  "
$ws.Range("D24").Value = "'2.03"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +3.79%  "
$ws.Range("D25").Value = "'151.67"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +0.83%  "
$ws.Range("D26").Value = "'15.13"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  -0.84%  "
$ws.Range("D27").Value = "'6.52"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  -1.36%  "
$ws.Range("D28").Value = "'0.105"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  -1.68%  "
$ws.Range("E29").Value = "  +1.12%  "
$ws.Range("E30").Value = "  -1.17%  "
$ws.Range("D31").Value = "'0.0467"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  -1.21%  "
$ws.Range("D32").Value = "'3.23"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  -0.72%  "
$ws.Range("D33").Value = "'3.14"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +0.54%  "
$ws.Range("D34").Value = "1.393.01"
$ws.Range("E34").Value = "  -4.68%  "
$ws.Range("D35").Value = "'1.57"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "  -2.41%  "
$ws.Range("E36").Value = "  -6.87%  "
$ws.Range("D37").Value = "'2.36"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  +1.50%  "
$ws.Range("D38").Value = "'2.61"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +10.06%  "
$ws.Range("D39").Value = "'0.0166"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  -1.13%  "
$ws.Range("D40").Value = "'0.536"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  -0.94%  "
$ws.Range("D41").Value = "'0.806"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  -1.07%  "
$ws.Range("E42").Value = "  +1.08%  "
$ws.Range("D43").Value = "'5.59"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  -1.56%  "
$ws.Range("E44").Value = "  +0.39%  "
$ws.Range("D45").Value = "'0.979"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +0.57%  "
$ws.Range("D46").Value = "'62.65"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "  -2.42%  "
$ws.Range("D47").Value = "1.719.57"
$ws.Range("E47").Value = "  +0.76%  "
$ws.Range("E48").Value = "  -0.18%  "
$ws.Range("D49").Value = "'86.37"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  -0.57%  "
$ws.Range("D50").Value = "0.0₆0103"
$ws.Range("E50").Value = "  -1.06%  "
$ws.Range("E51").Value = "  -0.79%  "
